$wb = $excel.ActiveWorkbook

# Rename the *img sheets to img* (himg->imgh, timg->imgt, simg->imgs,
# gimg->imgg, wimg->imgw, bimg->imgb, eimg->imge)
$wb.Worksheets.Item("himg").Name = "imgh"
$wb.Worksheets.Item("timg").Name = "imgt"
$wb.Worksheets.Item("simg").Name = "imgs"
$wb.Worksheets.Item("gimg").Name = "imgg"
$wb.Worksheets.Item("wimg").Name = "imgw"
$wb.Worksheets.Item("bimg").Name = "imgb"
$wb.Worksheets.Item("eimg").Name = "imge"

# Move the active tab from "holiday" to the last sheet, "imge"
$wb.Worksheets.Item("imge").Activate()
